# Update the LinkedIn carousel draft: the deck's topic changes from the
# "new solar capacity" announcement to the NTPC Green Energy / GAIL joint
# venture announcement. Each slide is a single TextBox with a title
# paragraph and two body paragraphs; we rewrite each run's text in place so
# existing paragraph/run formatting (defRPr, sizes, colors) is preserved.

$p = $ppt.ActivePresentation

function Set-SlideCopy {
    param(
        $Slide,
        [string]$Title,
        [string]$Body1,
        [string]$Body2
    )

    $shape = $Slide.Shapes.Item(1)

    # Remember the shape's current size so the auto-fit recalculation that
    # happens while we rewrite the text doesn't leave a different size
    # behind once we're done.
    $origWidth = $shape.Width
    $origHeight = $shape.Height

    $tr = $shape.TextFrame.TextRange
    $tr.Paragraphs(1).Runs(1).Text = $Title
    $tr.Paragraphs(2).Runs(1).Text = $Body1
    $tr.Paragraphs(3).Runs(1).Text = $Body2

    $shape.Width = $origWidth
    $shape.Height = $origHeight
}

Set-SlideCopy $p.Slides.Item(1) `
    "Joint Venture Formation" `
    "NTPC Green Energy and GAIL have established a 50:50 joint venture." `
    "The joint venture is focused on renewable energy projects."

Set-SlideCopy $p.Slides.Item(2) `
    "Participants" `
    "NTPC Green Energy is a subsidiary of NTPC Limited." `
    "GAIL is India's largest state-owned natural gas company."

Set-SlideCopy $p.Slides.Item(3) `
    "Renewable Energy Focus" `
    "The joint venture aims to develop renewable energy projects." `
    "This includes solar, wind, and other clean energy initiatives."

Set-SlideCopy $p.Slides.Item(4) `
    "Project Objectives" `
    "The projects are intended to enhance India's renewable energy capacity." `
    "They will contribute to the country's sustainability goals."

Set-SlideCopy $p.Slides.Item(5) `
    "Investment Plans" `
    "Details on specific investment amounts have not been disclosed." `
    "The joint venture will leverage resources from both companies."

Set-SlideCopy $p.Slides.Item(6) `
    "Timeline and Future Outlook" `
    "No specific timeline for project commencement has been mentioned." `
    "The joint venture aims to position itself strategically in the renewable sector."
